$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values per row.
# Price cells whose new text would otherwise look like a plain number
# (e.g. "155.72") are pre-formatted as Text so Excel keeps them as
# strings, matching the original inlineStr cell content/type.
$ws.Range("D2").Value = "64.274.47"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "3.397.09"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.14"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.72"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +8.93%  "
$ws.Range("D9").Value = "3.398.33"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "3.978.86"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000187"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.37"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").Value = "64.256.44"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "3.348.34"
$ws.Range("E18").Value = "  -2.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.28"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.06"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.00"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.542"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.69"
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("E27").Value = "  +7.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.176"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.11"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.00"
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.09"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  +6.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.19"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0757"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.75"
$ws.Range("E39").Value = "  +3.52%  "
$ws.Range("D40").Value = "2.862.44"
$ws.Range("E40").Value = "  -5.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.19"
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.58"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.63"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0314"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.82"
$ws.Range("E45").Value = "  +4.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.764"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "320.72"
$ws.Range("E47").Value = "  +4.83%  "
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("E51").Value = "  -1.34%  "
